$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.963.56"
$ws.Range("E2").Value = "  -2.01%  "

$ws.Range("D3").Value = "3.096.87"
$ws.Range("E3").Value = "  -0.32%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "526.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.04%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "3.095.86"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.445"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.24%  "

$ws.Range("E11").Value = "  -1.51%  "

$ws.Range("E12").Value = "  +1.89%  "

$ws.Range("D13").Value = "3.632.28"
$ws.Range("E13").Value = "  -0.19%  "

$ws.Range("E14").Value = "  +3.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.65%  "

$ws.Range("E16").Value = "  -1.45%  "

$ws.Range("D17").Value = "58.002.33"
$ws.Range("E17").Value = "  -1.85%  "

$ws.Range("D18").Value = "3.101.32"
$ws.Range("E18").Value = "  -0.06%  "

$ws.Range("E20").Value = "  -2.78%  "

$ws.Range("E21").Value = "  -2.81%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "342.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.45%  "

$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.513"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.61%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.28%  "

$ws.Range("E26").Value = "  -0.44%  "

$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").Value = "0.0₃0923"
$ws.Range("E28").Value = "  -1.45%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.06%  "

$ws.Range("E31").Value = "  +0.06%  "

$ws.Range("E32").Value = "  +1.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.69%  "

$ws.Range("E34").Value = "  -3.27%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.97"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.35%  "

$ws.Range("E36").Value = "  -0.67%  "

$ws.Range("E37").Value = "  -1.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.18%  "

$ws.Range("E39").Value = "  -4.98%  "

$ws.Range("E40").Value = "  -3.05%  "

$ws.Range("E41").Value = "  +7.65%  "

$ws.Range("E42").Value = "  +1.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.686"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.21%  "

$ws.Range("D44").Value = "3.140.22"
$ws.Range("E44").Value = "  -0.26%  "

$ws.Range("E45").Value = "  +0.27%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0263"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.62%  "

$ws.Range("D48").Value = "2.268.05"
$ws.Range("E48").Value = "  -1.11%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.992"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.21%  "
